# Wireframe caso de uso 3: add "Reservar cita" as its own run right
# after the existing "CU003: " label, matching the style already used
# for "CU002: " + "Ver psicólogos".

$d = $word.ActiveDocument

# Locate the "CU003: " label text and collapse the found range to its
# end so subsequent insertion happens right after it, inside the same
# paragraph.
$rng = $d.Content
$found = $rng.Find.Execute("CU003: ", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'CU003: ' label"
}

$rng.Collapse(0)

# Insert the new run's text right after "CU003: ".
$rng.InsertAfter("Reservar cita")

# Toggling a character attribute forces the newly inserted text to be
# kept as its own <w:r> run (matching the source document's pattern of
# separate runs for the label and its description) instead of being
# silently coalesced back into the preceding "CU003: " run, even though
# the final formatting is identical to it (bCs / sz 24 / szCs 24, no
# bold).
$rng.Font.Bold = $true
$rng.Font.Bold = $false

Write-Output "Inserted 'Reservar cita' after CU003 label"
